$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.777.87'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.732.26'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.08%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.596'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('E10').Value = '  +4.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.62'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.379'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.219.25'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.96'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.670.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000149'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.743.64'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '354.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.519'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.169'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0907'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('E29').Value = '  +3.31%  '
$ws.Range('E30').Value = '  +11.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.91'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.08'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('E35').Value = '  +2.66%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.80'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.973'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '343.73'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.26'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.67'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.77'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0584'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.630'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0250'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0997'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '132.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.02%  '
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.04%  '
